# Apply updated probability matrix values to Sheet1 ("team matrices from games pulled march 7")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B2").Value = 0.1903114186851211
    $ws.Range("C2").Value = 0.5847750865051903
    $ws.Range("J2").Value = 0.01038062283737024
    $ws.Range("P2").Value = 0.1245674740484429
    $ws.Range("S2").Value = 0.08996539792387544
    $ws.Range("B3").Value = 0.005681818181818182
    $ws.Range("C3").Value = 0.02840909090909091
    $ws.Range("J3").Value = 0.03409090909090909
    $ws.Range("P3").Value = 0.7102272727272727
    $ws.Range("S3").Value = 0.2215909090909091
    $ws.Range("J4").Value = 0.04761904761904762
    $ws.Range("P4").Value = 0.5714285714285714
    $ws.Range("S4").Value = 0.3809523809523809
    $ws.Range("B6").Value = 0.04979253112033195
    $ws.Range("D6").Value = 0.004149377593360996
    $ws.Range("E6").Value = 0.004149377593360996
    $ws.Range("F6").Value = 0.06639004149377593
    $ws.Range("J6").Value = 0.2282157676348548
    $ws.Range("O6").Value = 0.03319502074688797
    $ws.Range("Q6").Value = 0.1701244813278008
    $ws.Range("R6").Value = 0.0912863070539419
    $ws.Range("S6").Value = 0.3526970954356847
    $ws.Range("B7").Value = 0.08629441624365482
    $ws.Range("D7").Value = 0.02538071065989848
    $ws.Range("F7").Value = 0.05076142131979695
    $ws.Range("J7").Value = 0.1472081218274112
    $ws.Range("O7").Value = 0.02538071065989848
    $ws.Range("Q7").Value = 0.1624365482233502
    $ws.Range("R7").Value = 0.06598984771573604
    $ws.Range("S7").Value = 0.4365482233502538
    $ws.Range("B8").Value = 0.08836206896551724
    $ws.Range("D8").Value = 0.02370689655172414
    $ws.Range("F8").Value = 0.03232758620689655
    $ws.Range("J8").Value = 0.125
    $ws.Range("O8").Value = 0.03017241379310345
    $ws.Range("Q8").Value = 0.1918103448275862
    $ws.Range("R8").Value = 0.09267241379310345
    $ws.Range("S8").Value = 0.415948275862069
    $ws.Range("B9").Value = 0.08677685950413223
    $ws.Range("D9").Value = 0.02479338842975207
    $ws.Range("F9").Value = 0.07024793388429752
    $ws.Range("J9").Value = 0.1570247933884298
    $ws.Range("O9").Value = 0.01239669421487603
    $ws.Range("Q9").Value = 0.2520661157024793
    $ws.Range("R9").Value = 0.07024793388429752
    $ws.Range("S9").Value = 0.3264462809917356
    $ws.Range("B10").Value = 0.09570724841660802
    $ws.Range("D10").Value = 0.01477832512315271
    $ws.Range("E10").Value = 0.001407459535538353
    $ws.Range("F10").Value = 0.08163265306122448
    $ws.Range("J10").Value = 0.1358198451794511
    $ws.Range("O10").Value = 0.01759324419422941
    $ws.Range("Q10").Value = 0.1977480647431386
    $ws.Range("R10").Value = 0.07459535538353272
    $ws.Range("S10").Value = 0.3807178043631246
    $ws.Range("G11").Value = 0.1346153846153846
    $ws.Range("J11").Value = 0.108974358974359
    $ws.Range("K11").Value = 0.2083333333333333
    $ws.Range("L11").Value = 0.5224358974358975
    $ws.Range("S11").Value = 0.02564102564102564
    $ws.Range("G12").Value = 0.6949152542372882
    $ws.Range("J12").Value = 0.192090395480226
    $ws.Range("K12").Value = 0.02259887005649718
    $ws.Range("L12").Value = 0.05084745762711865
    $ws.Range("S12").Value = 0.03954802259887006
    $ws.Range("G14").Value = 0.5
    $ws.Range("J14").Value = 0.5
    $ws.Range("F15").Value = 0.01538461538461539
    $ws.Range("H15").Value = 0.1192307692307692
    $ws.Range("I15").Value = 0.08076923076923077
    $ws.Range("J15").Value = 0.3692307692307693
    $ws.Range("K15").Value = 0.0576923076923077
    $ws.Range("M15").Value = 0.01153846153846154
    $ws.Range("N15").Value = 0.003846153846153846
    $ws.Range("O15").Value = 0.07307692307692308
    $ws.Range("S15").Value = 0.2692307692307692
    $ws.Range("F16").Value = 0.01714285714285714
    $ws.Range("H16").Value = 0.1714285714285714
    $ws.Range("I16").Value = 0.05714285714285714
    $ws.Range("J16").Value = 0.4685714285714286
    $ws.Range("K16").Value = 0.1085714285714286
    $ws.Range("M16").Value = 0.04
    $ws.Range("O16").Value = 0.06285714285714286
    $ws.Range("S16").Value = 0.07428571428571429
    $ws.Range("F17").Value = 0.01006036217303823
    $ws.Range("H17").Value = 0.1971830985915493
    $ws.Range("I17").Value = 0.1146881287726358
    $ws.Range("J17").Value = 0.4064386317907445
    $ws.Range("K17").Value = 0.09456740442655935
    $ws.Range("M17").Value = 0.01207243460764587
    $ws.Range("N17").Value = 0.004024144869215292
    $ws.Range("O17").Value = 0.05835010060362173
    $ws.Range("S17").Value = 0.1026156941649899
    $ws.Range("F18").Value = 0.01515151515151515
    $ws.Range("H18").Value = 0.1464646464646465
    $ws.Range("I18").Value = 0.09595959595959595
    $ws.Range("J18").Value = 0.4444444444444444
    $ws.Range("K18").Value = 0.1111111111111111
    $ws.Range("M18").Value = 0.0101010101010101
    $ws.Range("O18").Value = 0.08585858585858586
    $ws.Range("S18").Value = 0.09090909090909091
    $ws.Range("F19").Value = 0.01745454545454546
    $ws.Range("H19").Value = 0.1949090909090909
    $ws.Range("I19").Value = 0.1003636363636364
    $ws.Range("J19").Value = 0.3774545454545454
    $ws.Range("K19").Value = 0.1003636363636364
    $ws.Range("M19").Value = 0.02181818181818182
    $ws.Range("N19").Value = 0.001454545454545454
    $ws.Range("O19").Value = 0.07636363636363637
    $ws.Range("S19").Value = 0.1098181818181818

